$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.070632
$ws.Range("H2").Value = 30.211896
$ws.Range("I2").Value = 0.07634150103324112
$ws.Range("J2").Value = 0.08026042296304617
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 197.2278263333334
$ws.Range("N2").Value = 591.683479
$ws.Range("O2").Value = 0.6783778564662776
$ws.Range("P2").Value = 0.6850369527608899
$ws.Range("Q2").Value = 1986.208859162909
$ws.Range("R2").Value = 17875.87973246619
$ws.Range("S2").Value = 0.05178838383034823
$ws.Range("T2").Value = 0.0549813555739053

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.070632
$ws.Range("H3").Value = 30.211896
$ws.Range("I3").Value = 0.07634150103324112
$ws.Range("J3").Value = 0.08026042296304617
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.896351
$ws.Range("N3").Value = 2.689053
$ws.Range("O3").Value = 0.003083057200020643
$ws.Range("P3").Value = 0.003113321122377543
$ws.Range("Q3").Value = 9.026821063832
$ws.Range("R3").Value = 81.241389574488
$ws.Range("S3").Value = 0.0002353652144209174
$ws.Range("T3").Value = 0.0002498764701018072

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.070632
$ws.Range("H4").Value = 30.211896
$ws.Range("I4").Value = 0.07634150103324112
$ws.Range("J4").Value = 0.08026042296304617
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 42.82536200000001
$ws.Range("N4").Value = 128.476086
$ws.Range("O4").Value = 0.1473006006102414
$ws.Range("P4").Value = 0.1487465335432934
$ws.Range("Q4").Value = 431.278460968784
$ws.Range("R4").Value = 3881.506148719056
$ws.Range("S4").Value = 0.01124514895368378
$ws.Range("T4").Value = 0.01193845969647166

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.070632
$ws.Range("H5").Value = 30.211896
$ws.Range("I5").Value = 0.07634150103324112
$ws.Range("J5").Value = 0.08026042296304617
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 41.30642933333333
$ws.Range("N5").Value = 123.919288
$ws.Range("O5").Value = 0.1420761335272424
$ws.Range("P5").Value = 0.1434707820189434
$ws.Range("Q5").Value = 415.9818490500053
$ws.Range("R5").Value = 3743.836641450048
$ws.Range("S5").Value = 0.01084630529446888
$ws.Range("T5").Value = 0.0115150256476794

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.070632
$ws.Range("H6").Value = 30.211896
$ws.Range("I6").Value = 0.07634150103324112
$ws.Range("J6").Value = 0.08026042296304617
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 8.478501
$ws.Range("N6").Value = 16.957002
$ws.Range("O6").Value = 0.02916235219621802
$ws.Range("P6").Value = 0.01963241055449567
$ws.Range("Q6").Value = 85.383863482632
$ws.Range("R6").Value = 512.3031808957919
$ws.Range("S6").Value = 0.00222629774031932
$ws.Range("T6").Value = 0.001575705574887994

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 42.14988333333334
$ws.Range("H7").Value = 126.44965
$ws.Range("I7").Value = 0.3195216905992255
$ws.Range("J7").Value = 0.3359240476840365
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 197.2278263333334
$ws.Range("N7").Value = 591.683479
$ws.Range("O7").Value = 0.6783778564662776
$ws.Range("P7").Value = 0.6850369527608899
$ws.Range("Q7").Value = 8313.12987003693
$ws.Range("R7").Value = 74818.16883033236
$ws.Range("S7").Value = 0.2167564395631838
$ws.Range("T7").Value = 0.2301203859845762

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 42.14988333333334
$ws.Range("H8").Value = 126.44965
$ws.Range("I8").Value = 0.3195216905992255
$ws.Range("J8").Value = 0.3359240476840365
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.896351
$ws.Range("N8").Value = 2.689053
$ws.Range("O8").Value = 0.003083057200020643
$ws.Range("P8").Value = 0.003113321122377543
$ws.Range("Q8").Value = 37.78109007571667
$ws.Range("R8").Value = 340.02981068145
$ws.Range("S8").Value = 0.0009851036487647106
$ws.Range("T8").Value = 0.001045839433169272

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 42.14988333333334
$ws.Range("H9").Value = 126.44965
$ws.Range("I9").Value = 0.3195216905992255
$ws.Range("J9").Value = 0.3359240476840365
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 42.82536200000001
$ws.Range("N9").Value = 128.476086
$ws.Range("O9").Value = 0.1473006006102414
$ws.Range("P9").Value = 0.1487465335432934
$ws.Range("Q9").Value = 1805.084012007767
$ws.Range("R9").Value = 16245.7561080699
$ws.Range("S9").Value = 0.04706573693326563
$ws.Range("T9").Value = 0.04996753762683243

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 42.14988333333334
$ws.Range("H10").Value = 126.44965
$ws.Range("I10").Value = 0.3195216905992255
$ws.Range("J10").Value = 0.3359240476840365
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 41.30642933333333
$ws.Range("N10").Value = 123.919288
$ws.Range("O10").Value = 0.1420761335272424
$ws.Range("P10").Value = 0.1434707820189434
$ws.Range("Q10").Value = 1741.061177316578
$ws.Range("R10").Value = 15669.5505958492
$ws.Range("S10").Value = 0.04539640637842578
$ws.Range("T10").Value = 0.04819528582019757

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 42.14988333333334
$ws.Range("H11").Value = 126.44965
$ws.Range("I11").Value = 0.3195216905992255
$ws.Range("J11").Value = 0.3359240476840365
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 8.478501
$ws.Range("N11").Value = 16.957002
$ws.Range("O11").Value = 0.02916235219621802
$ws.Range("P11").Value = 0.01963241055449567
$ws.Range("Q11").Value = 357.36782799155
$ws.Range("R11").Value = 2144.2069679493
$ws.Range("S11").Value = 0.009318004075585619
$ws.Range("T11").Value = 0.006594998819260984

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 26.941003
$ws.Range("H12").Value = 80.823009
$ws.Range("I12").Value = 0.2042291495073052
$ws.Range("J12").Value = 0.2147130682392819
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 197.2278263333334
$ws.Range("N12").Value = 591.683479
$ws.Range("O12").Value = 0.6783778564662776
$ws.Range("P12").Value = 0.6850369527608899
$ws.Range("Q12").Value = 5313.515460929812
$ws.Range("R12").Value = 47821.63914836831
$ws.Range("S12").Value = 0.1385445326706966
$ws.Range("T12").Value = 0.1470863859845787

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 26.941003
$ws.Range("H13").Value = 80.823009
$ws.Range("I13").Value = 0.2042291495073052
$ws.Range("J13").Value = 0.2147130682392819
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.896351
$ws.Range("N13").Value = 2.689053
$ws.Range("O13").Value = 0.003083057200020643
$ws.Range("P13").Value = 0.003113321122377543
$ws.Range("Q13").Value = 24.148594980053
$ws.Range("R13").Value = 217.337354820477
$ws.Range("S13").Value = 0.0006296501498425898
$ws.Range("T13").Value = 0.000668470730599847

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 26.941003
$ws.Range("H14").Value = 80.823009
$ws.Range("I14").Value = 0.2042291495073052
$ws.Range("J14").Value = 0.2147130682392819
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 42.82536200000001
$ws.Range("N14").Value = 128.476086
$ws.Range("O14").Value = 0.1473006006102414
$ws.Range("P14").Value = 0.1487465335432934
$ws.Range("Q14").Value = 1153.758206118086
$ws.Range("R14").Value = 10383.82385506277
$ws.Range("S14").Value = 0.03008307638454484
$ws.Range("T14").Value = 0.03193782460703778

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 26.941003
$ws.Range("H15").Value = 80.823009
$ws.Range("I15").Value = 0.2042291495073052
$ws.Range("J15").Value = 0.2147130682392819
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 41.30642933333333
$ws.Range("N15").Value = 123.919288
$ws.Range("O15").Value = 0.1420761335272424
$ws.Range("P15").Value = 0.1434707820189434
$ws.Range("Q15").Value = 1112.836636588621
$ws.Range("R15").Value = 10015.52972929759
$ws.Range("S15").Value = 0.02901608791555504
$ws.Range("T15").Value = 0.03080505180997654

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 26.941003
$ws.Range("H16").Value = 80.823009
$ws.Range("I16").Value = 0.2042291495073052
$ws.Range("J16").Value = 0.2147130682392819
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 8.478501
$ws.Range("N16").Value = 16.957002
$ws.Range("O16").Value = 0.02916235219621802
$ws.Range("P16").Value = 0.01963241055449567
$ws.Range("Q16").Value = 228.419320876503
$ws.Range("R16").Value = 1370.515925259018
$ws.Range("S16").Value = 0.005955802386666101
$ws.Range("T16").Value = 0.004215335107089026

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 33.430692
$ws.Range("H17").Value = 100.292076
$ws.Range("I17").Value = 0.2534249298216801
$ws.Range("J17").Value = 0.2664342694547198
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 197.2278263333334
$ws.Range("N17").Value = 591.683479
$ws.Range("O17").Value = 0.6783778564662776
$ws.Range("P17").Value = 0.6850369527608899
$ws.Range("Q17").Value = 6593.462715979157
$ws.Range("R17").Value = 59341.16444381241
$ws.Range("S17").Value = 0.1719178606675482
$ws.Range("T17").Value = 0.1825173200583351

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 33.430692
$ws.Range("H18").Value = 100.292076
$ws.Range("I18").Value = 0.2534249298216801
$ws.Range("J18").Value = 0.2664342694547198
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 0.896351
$ws.Range("N18").Value = 2.689053
$ws.Range("O18").Value = 0.003083057200020643
$ws.Range("P18").Value = 0.003113321122377543
$ws.Range("Q18").Value = 29.965634204892
$ws.Range("R18").Value = 269.690707844028
$ws.Range("S18").Value = 0.000781323554551457
$ws.Range("T18").Value = 0.0008294954388186089

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 33.430692
$ws.Range("H19").Value = 100.292076
$ws.Range("I19").Value = 0.2534249298216801
$ws.Range("J19").Value = 0.2664342694547198
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 42.82536200000001
$ws.Range("N19").Value = 128.476086
$ws.Range("O19").Value = 0.1473006006102414
$ws.Range("P19").Value = 0.1487465335432934
$ws.Range("Q19").Value = 1431.681486810504
$ws.Range("R19").Value = 12885.13338129454
$ws.Range("S19").Value = 0.03732964437234174
$ws.Range("T19").Value = 0.03963117399852936

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 33.430692
$ws.Range("H20").Value = 100.292076
$ws.Range("I20").Value = 0.2534249298216801
$ws.Range("J20").Value = 0.2664342694547198
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 41.30642933333333
$ws.Range("N20").Value = 123.919288
$ws.Range("O20").Value = 0.1420761335272424
$ws.Range("P20").Value = 0.1434707820189434
$ws.Range("Q20").Value = 1380.902516662432
$ws.Range("R20").Value = 12428.12264996189
$ws.Range("S20").Value = 0.03600563416847704
$ws.Range("T20").Value = 0.03822553299531455

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 33.430692
$ws.Range("H21").Value = 100.292076
$ws.Range("I21").Value = 0.2534249298216801
$ws.Range("J21").Value = 0.2664342694547198
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 8.478501
$ws.Range("N21").Value = 16.957002
$ws.Range("O21").Value = 0.02916235219621802
$ws.Range("P21").Value = 0.01963241055449567
$ws.Range("Q21").Value = 283.442155552692
$ws.Range("R21").Value = 1700.652933316152
$ws.Range("S21").Value = 0.00739046705876167
$ws.Range("T21").Value = 0.005230746963722184

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 19.3233515
$ws.Range("H22").Value = 38.646703
$ws.Range("I22").Value = 0.1464827290385481
$ws.Range("J22").Value = 0.1026681916589156
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 197.2278263333334
$ws.Range("N22").Value = 591.683479
$ws.Range("O22").Value = 0.6783778564662776
$ws.Range("P22").Value = 0.6850369527608899
$ws.Range("Q22").Value = 3811.102613819957
$ws.Range("R22").Value = 22866.61568291974
$ws.Range("S22").Value = 0.0993706397345008
$ws.Range("T22").Value = 0.07033150515949456

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 19.3233515
$ws.Range("H23").Value = 38.646703
$ws.Range("I23").Value = 0.1464827290385481
$ws.Range("J23").Value = 0.1026681916589156
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.896351
$ws.Range("N23").Value = 2.689053
$ws.Range("O23").Value = 0.003083057200020643
$ws.Range("P23").Value = 0.003113321122377543
$ws.Range("Q23").Value = 17.3205054403765
$ws.Range("R23").Value = 103.923032642259
$ws.Range("S23").Value = 0.0004516146324409686
$ws.Range("T23").Value = 0.0003196390496880078

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 19.3233515
$ws.Range("H24").Value = 38.646703
$ws.Range("I24").Value = 0.1464827290385481
$ws.Range("J24").Value = 0.1026681916589156
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 42.82536200000001
$ws.Range("N24").Value = 128.476086
$ws.Range("O24").Value = 0.1473006006102414
$ws.Range("P24").Value = 0.1487465335432934
$ws.Range("Q24").Value = 827.5295230407431
$ws.Range("R24").Value = 4965.177138244459
$ws.Range("S24").Value = 0.02157699396640538
$ws.Range("T24").Value = 0.01527153761442216

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 19.3233515
$ws.Range("H25").Value = 38.646703
$ws.Range("I25").Value = 0.1464827290385481
$ws.Range("J25").Value = 0.1026681916589156
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 41.30642933333333
$ws.Range("N25").Value = 123.919288
$ws.Range("O25").Value = 0.1420761335272424
$ws.Range("P25").Value = 0.1434707820189434
$ws.Range("Q25").Value = 798.1786532179108
$ws.Range("R25").Value = 4789.071919307464
$ws.Range("S25").Value = 0.02081169977031562
$ws.Range("T25").Value = 0.01472988574577539

$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 19.3233515
$ws.Range("H26").Value = 38.646703
$ws.Range("I26").Value = 0.1464827290385481
$ws.Range("J26").Value = 0.1026681916589156
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 8.478501
$ws.Range("N26").Value = 16.957002
$ws.Range("O26").Value = 0.02916235219621802
$ws.Range("P26").Value = 0.01963241055449567
$ws.Range("Q26").Value = 163.8330550161015
$ws.Range("R26").Value = 655.332220064406
$ws.Range("S26").Value = 0.004271780934885312
$ws.Range("T26").Value = 0.002015624089535479
